# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig@8e4a450c507ef6f746e072652acbb72e9504f19a
# Regenerated "Metadata" sheet values (version bump, date bump, publisher, jurisdiction)
# and "Elements" sheet Short/Definition text for the root Extension row.

$wb = $excel.ActiveWorkbook
$metadata = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------------

# Version: 5.0.0 -> 6.0.0
$metadata.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$metadata.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank -> "Alvearie Team"
$metadata.Range("B9").Value = "Alvearie Team"

# Row 10 used to be a stray duplicate "Contact" / "No display for ContactDetail"
# row; turn it into the real "Jurisdiction" row.
$metadata.Range("A10").Value = "Jurisdiction"
$metadata.Range("B10").Value = "United States of America"

# Row 11 was the second (duplicate) "Contact" row - remove it entirely so
# everything below shifts up by one row.
$metadata.Rows.Item(11).Delete()

# --- Elements sheet --------------------------------------------------------

# Root Extension row's Short/Definition columns (K2/L2) get the real
# short description & definition instead of the generic placeholders.
$elements.Range("K2").Value = "Communication Language"
$elements.Range("L2").Value = "Language used for communication messaging content"
